$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 78 - this shifts the existing rows 78:185
# down to 79:186 (and all their formatting/values move with them).
$ws.Rows(78).Insert()

# Populate the newly inserted (blank) row 78 with the new record.
$ws.Range("A78").Value = 5
$ws.Range("B78").Value = "Macroferia Regional de Talca"
$ws.Range("C78").Value = "Maule"
$ws.Range("D78").Value = 44482
$ws.Range("E78").Value = 7
$ws.Range("F78").Value = 100114014
$ws.Range("G78").Value = "Betarraga"
$ws.Range("H78").Value = "Sin especificar"
$ws.Range("I78").Value = "Primera"
$ws.Range("J78").Value = 4000
$ws.Range("K78").Value = 700
$ws.Range("L78").Value = 700
$ws.Range("M78").Value = 700
$ws.Range("N78").Value = "$/paquete 5 unidades"
$ws.Range("O78").Value = "Región del Maule"
$ws.Range("P78").Value = 140
$ws.Range("Q78").Value = 5
$ws.Range("R78").Value = "Hortaliza"
